$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a map of cell address -> new text value, preserving the original
# "inline string" text representation (e.g. trailing zeros, percent signs).
# NumberFormat is forced to Text ("@") before assignment so Excel does not
# auto-convert the strings to numbers/percentages and strip formatting.
$updates = @{
    'D2' = '330.32'
    'E2' = '4.51%'
    'G2' = '21'
    'D3' = '40.58'
    'E3' = '8.25%'
    'G3' = '21'
    'D4' = '5.735'
    'E4' = '11.32%'
    'G4' = '21'
    'D5' = '0.08102'
    'E5' = '1.72%'
    'G5' = '21'
    'E6' = '3.24%'
    'G6' = '21'
    'D7' = '8.774'
    'E7' = '3.95%'
    'G7' = '21'
    'D8' = '1.967'
    'E8' = '2.83%'
    'G8' = '21'
    'D9' = '2.941'
    'E9' = '-0.42%'
    'G9' = '21'
    'D10' = '0.9466'
    'E10' = '0.72%'
    'G10' = '21'
    'D11' = '0.1299'
    'E11' = '2.56%'
    'G11' = '21'
    'D12' = '0.1994'
    'E12' = '3.48%'
    'G12' = '21'
    'D13' = '8.956'
    'E13' = '37.49%'
    'G13' = '21'
    'D14' = '0.09412'
    'E14' = '4.94%'
    'G14' = '21'
    'D15' = '0.03521'
    'E15' = '3.99%'
    'G15' = '21'
    'D16' = '0.09634'
    'G16' = '21'
    'D17' = '0.001326'
    'E17' = '-4.67%'
    'G17' = '21'
    'D18' = '0.006114'
    'E18' = '0.92%'
    'G18' = '21'
    'D19' = '3.369'
    'E19' = '-0.75%'
    'G19' = '21'
    'E20' = '1.50%'
    'G20' = '21'
    'D21' = '0.1420'
    'E21' = '9.07%'
    'G21' = '21'
    'D22' = '0.2409'
    'E22' = '4.64%'
    'G22' = '21'
    'D23' = '0.04415'
    'E23' = '1.59%'
    'G23' = '21'
    'D24' = '0.001257'
    'E24' = '4.98%'
    'G24' = '21'
    'D25' = '0.004370'
    'E25' = '-0.75%'
    'G25' = '21'
    'D26' = '0.0001090'
    'G26' = '21'
    'D27' = '0.0003993'
    'E27' = '0.83%'
    'G27' = '21'
    'G28' = '21'
    'G29' = '21'
    'G30' = '21'
    'G31' = '21'
    'G32' = '21'
    'G33' = '21'
    'G34' = '21'
    'G35' = '21'
    'G36' = '21'
    'G37' = '21'
    'G38' = '21'
    'D39' = '0.02464'
    'E39' = '5.95%'
    'G39' = '21'
    'D40' = '0.05325'
    'E40' = '3.12%'
    'G40' = '21'
    'D41' = '0.007491'
    'E41' = '0.29%'
    'G41' = '21'
    'D42' = '0.1437'
    'E42' = '3.12%'
    'G42' = '21'
    'D43' = '0.008843'
    'E43' = '3.13%'
    'G43' = '21'
    'D44' = '0.002121'
    'E44' = '6.48%'
    'G44' = '21'
    'D45' = '0.01043'
    'E45' = '31.57%'
    'G45' = '21'
    'D46' = '0.00006877'
    'E46' = '7.95%'
    'G46' = '21'
    'D47' = '0.00000000750'
    'E47' = '0.58%'
    'G47' = '21'
    'D48' = '0.003502'
    'E48' = '22.79%'
    'G48' = '21'
    'D49' = '0.001701'
    'E49' = '1.22%'
    'G49' = '21'
    'D50' = '0.00002101'
    'E50' = '0.58%'
    'G50' = '21'
    'D51' = '0.0002001'
    'E51' = '0.58%'
    'G51' = '21'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
